# Error Calculations and Plots
# Apply edits to the missing_data worksheet:
#  - Remove two rows that no longer belong in the data set ("RM 232" and "SC 92")
#  - Fill in / clear a handful of cells across columns E and F to reflect
#    the corrected / re-checked values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "RM 232" row entirely (was row 26); everything below shifts up.
$ws.Rows(26).Delete()

# Drop the "SC 92" row entirely (was row 28, now row 27 after the first delete).
$ws.Rows(27).Delete()

# --- Column E / F value corrections on the remaining rows ---

# RM 8 (row 3): column E value filled in
$ws.Range("E3").Value = -5.7

# RM 9 (row 4): column F value removed
$ws.Range("F4").ClearContents()

# RM 14 (row 5): column E value removed
$ws.Range("E5").ClearContents()

# RM 42 (row 9): column F value filled in
$ws.Range("F9").Value = 17.26

# RM 52 a (row 10): column F value filled in
$ws.Range("F10").Value = 16.43

# RM 58 (row 11): column F value filled in
$ws.Range("F11").Value = 17.65

# RM 81 (row 12): column F value filled in
$ws.Range("F12").Value = 17.45

# RM 95 (row 15): column F value removed
$ws.Range("F15").ClearContents()

# RM 116 (row 17): column F value removed
$ws.Range("F17").ClearContents()

# RM 120 (row 18): column F value removed
$ws.Range("F18").ClearContents()

# RM 134 (row 20): column F value removed
$ws.Range("F20").ClearContents()

# RM 135 (row 21): column E value filled in
$ws.Range("E21").Value = -8.699999999999999

# RM 140 (row 23): column E value removed
$ws.Range("E23").ClearContents()

# SC 132 (row 31 after the two deletions): column F value filled in
$ws.Range("F31").Value = 17.18

# SC 193 (row 32 after the two deletions): column E and F values filled in
$ws.Range("E32").Value = -6.4
$ws.Range("F32").Value = 17.39

"Edits applied"
